# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first data row (the 6ae3d7f8-... file) on both the
# "zh-cn" and "de-de" worksheets, reflecting a newer handback report run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-12 18:38:57"
$zhcn.Range("H2").Value = "2016-03-12 18:39:14"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-12 18:39:00"
$dede.Range("H2").Value = "2016-03-12 18:39:20"
